# Refresh the cryptos price list (scraped values) in-place.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h)
# Numeric-looking Price strings are written with a leading apostrophe so
# Excel keeps them as text instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '41.513.02'
$ws.Range("E2").Value = '  +0.00%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.479.96'
$ws.Range("E3").Value = '  +0.72%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.12%  '

# Row 5: BNB
$ws.Range("D5").Value = '''310.99'
$ws.Range("E5").Value = '  +0.11%  '

# Row 6: Solana
$ws.Range("D6").Value = '''92.59'
$ws.Range("E6").Value = '  -2.17%  '

# Row 7: XRP
$ws.Range("D7").Value = '''0.537'
$ws.Range("E7").Value = '  -2.51%  '

# Row 8: USDC
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  +0.13%  '

# Row 9: Cardano
$ws.Range("D9").Value = '''0.491'
$ws.Range("E9").Value = '  -3.35%  '

# Row 10: Avalanche
$ws.Range("D10").Value = '''32.07'
$ws.Range("E10").Value = '  -4.83%  '

# Row 11: Dogecoin
$ws.Range("D11").Value = '''0.0774'
$ws.Range("E11").Value = '  -0.93%  '

# Row 12: TRON
$ws.Range("E12").Value = '  +1.31%  '

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = '2.847.79'
$ws.Range("E13").Value = '  +0.21%  '

# Row 14: Polkadot
$ws.Range("D14").Value = '''6.76'
$ws.Range("E14").Value = '  -2.85%  '

# Row 15: WrappedEther
$ws.Range("D15").Value = '2.506.66'
$ws.Range("E15").Value = '  +1.51%  '

# Row 16: Chainlink
$ws.Range("D16").Value = '''15.21'
$ws.Range("E16").Value = '  +4.65%  '

# Row 17: Polygon
$ws.Range("D17").Value = '''0.760'
$ws.Range("E17").Value = '  -3.42%  '

# Row 18: WrappedBTC
$ws.Range("D18").Value = '41.467.11'
$ws.Range("E18").Value = '  -0.09%  '

# Row 19: Uniswap
$ws.Range("D19").Value = '''6.26'
$ws.Range("E19").Value = '  -1.61%  '

# Row 20: ShibaInu
$ws.Range("D20").Value = '0.0₃0916'
$ws.Range("E20").Value = '  -0.17%  '

# Row 21: Litecoin
$ws.Range("D21").Value = '''70.43'
$ws.Range("E21").Value = '  +1.05%  '

# Row 22: InternetComputer(DFINITY)
$ws.Range("D22").Value = '''11.05'
$ws.Range("E22").Value = '  -4.26%  '

# Row 23: BitcoinCash
$ws.Range("D23").Value = '''234.08'
$ws.Range("E23").Value = '  -1.15%  '

# Row 24: PancakeSwap
$ws.Range("D24").Value = '''2.69'
$ws.Range("E24").Value = '  -3.03%  '

# Row 25: Dai
$ws.Range("E25").Value = '  -0.10%  '

# Row 26: ImmutableX
$ws.Range("D26").Value = '''1.87'
$ws.Range("E26").Value = '  -2.98%  '

# Row 27: EthereumClassic
$ws.Range("D27").Value = '''24.21'
$ws.Range("E27").Value = '  -1.87%  '

# Row 28: Toncoin
$ws.Range("D28").Value = '''2.24'
$ws.Range("E28").Value = '  +1.15%  '

# Row 29: Cosmos
$ws.Range("D29").Value = '''9.58'
$ws.Range("E29").Value = '  -1.55%  '

# Row 30: InjectiveProtocol
$ws.Range("D30").Value = '''36.13'
$ws.Range("E30").Value = '  -0.51%  '

# Row 31: Monero
$ws.Range("D31").Value = '''153.62'
$ws.Range("E31").Value = '  -0.17%  '

# Row 32: Filecoin
$ws.Range("D32").Value = '''5.35'
$ws.Range("E32").Value = '  -4.55%  '

# Row 33: WEMIXToken
$ws.Range("E33").Value = '  -2.46%  '

# Row 34: Celestia (was Hedera)
$ws.Range("B34").Value = 'Celestia'
$ws.Range("C34").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D34").Value = '''18.05'
$ws.Range("E34").Value = '  +4.29%  '

# Row 35: Hedera (was Celestia)
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '''0.0756'
$ws.Range("E35").Value = '  +0.28%  '

# Row 36: ApeXProtocol
$ws.Range("E36").Value = '  -2.08%  '

# Row 37: LidoDAOToken
$ws.Range("D37").Value = '''2.97'
$ws.Range("E37").Value = '  -1.57%  '

# Row 38: ARBITRUM
$ws.Range("D38").Value = '''1.82'
$ws.Range("E38").Value = '  -2.82%  '

# Row 39: Stellar
$ws.Range("E39").Value = '  -1.74%  '

# Row 40: Kaspa
$ws.Range("D40").Value = '''0.100'
$ws.Range("E40").Value = '  -4.39%  '

# Row 41: RenderToken
$ws.Range("D41").Value = '''4.06'
$ws.Range("E41").Value = '  +1.22%  '

# Row 42: EnergySwap (was FirstDigitalUSD)
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").Value = '''20.65'
$ws.Range("E42").Value = '  -3.10%  '

# Row 43: FirstDigitalUSD (was EnergySwap)
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").Value = '''1.01'
$ws.Range("E43").Value = '  +0.34%  '

# Row 44: Maker
$ws.Range("D44").Value = '1.945.50'
$ws.Range("E44").Value = '  -1.88%  '

# Row 45: VeChain
$ws.Range("D45").Value = '''0.0279'
$ws.Range("E45").Value = '  -2.13%  '

# Row 46: NEARProtocol
$ws.Range("D46").Value = '''2.93'
$ws.Range("E46").Value = '  -4.26%  '

# Row 47: FraxShare
$ws.Range("D47").Value = '''8.68'
$ws.Range("E47").Value = '  +0.00%  '

# Row 48: RocketPoolETH
$ws.Range("D48").Value = '2.714.23'
$ws.Range("E48").Value = '  +0.46%  '

# Row 49: Aave
$ws.Range("D49").Value = '''95.35'
$ws.Range("E49").Value = '  -2.27%  '

# Row 50: Algorand
$ws.Range("D50").Value = '''0.174'
$ws.Range("E50").Value = '  -3.63%  '

# Row 51: ordi
$ws.Range("D51").Value = '''66.22'
$ws.Range("E51").Value = '  -4.82%  '

